$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.981.95"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "3.148.97"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.139.65"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "3.665.81"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("E16").Value = "  +2.42%  "
$ws.Range("D17").Value = "63.970.28"
$ws.Range("D18").Value = "3.156.66"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "489.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.32%  "
$ws.Range("E25").Value = "  -4.05%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.53%  "
$ws.Range("E32").Value = "  -6.22%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D38").Value = "0.0₃0749"
$ws.Range("E38").Value = "  -5.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.48%  "
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "433.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.38%  "
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").Value = "2.931.88"
$ws.Range("E44").Value = "  +2.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.260"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.42%  "
$ws.Range("E46").Value = "  -6.31%  "
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.80%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "
